$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a few existing cells ---

# Row 35: Last Visited date changed to 45323
$ws.Range("E35").Value = 45323

# Row 38: First Visited date added
$ws.Range("D38").Value = 43405
$ws.Range("D38").NumberFormat = $ws.Range("E38").NumberFormat

# Row 40: Last Visited date changed to 45375
$ws.Range("E40").Value = 45375

# --- Append new rows 63-67 (USA / CT / MD entries) ---
# Values are entered in a specific order so new shared strings line up
# with how they were authored (cities first, then comments).

$ws.Range("C63").Value = "Mansfield"
$ws.Range("C64").Value = "Hartford"
$ws.Range("C65").Value = "South Windsor"
$ws.Range("G63").Value = "Visited University of Connecticut here! It's a very beautiful campus!"
$ws.Range("G65").Value = "Stayed at Cambria Hotel Manchester South Windsor here, it's a bougie hotel without the bougie price tag!"
$ws.Range("G64").Value = "Nice downtown area."
$ws.Range("C66").Value = "Port Deposit"
$ws.Range("G66").Value = "Stopped at Chesapeake Travel Plaza - bougie rest area!"
$ws.Range("C67").Value = "College Park"
$ws.Range("G67").Value = "Visited University of Maryland campus!"

# Fill in the remaining columns for the new rows
$ws.Range("A63").Value = "USA"
$ws.Range("A64").Value = "USA"
$ws.Range("A65").Value = "USA"
$ws.Range("A66").Value = "USA"
$ws.Range("A67").Value = "USA"

$ws.Range("B63").Value = "CT"
$ws.Range("B64").Value = "CT"
$ws.Range("B65").Value = "CT"
$ws.Range("B66").Value = "MD"
$ws.Range("B67").Value = "MD"

$ws.Range("D63").Value = 45352
$ws.Range("D64").Value = 45352
$ws.Range("D65").Value = 45352
$ws.Range("D66").Value = 45352
$ws.Range("D67").Value = 45323

$ws.Range("E63").Value = 45352
$ws.Range("E64").Value = 45352
$ws.Range("E65").Value = 45352
$ws.Range("E66").Value = 45352
$ws.Range("E67").Value = 45323

$ws.Range("F63").Value = "Visited"
$ws.Range("F64").Value = "Visited"
$ws.Range("F65").Value = "Visited"
$ws.Range("F66").Value = "Visited"
$ws.Range("F67").Value = "Visited"

$ws.Range("D63:E67").NumberFormat = $ws.Range("E61").NumberFormat
$ws.Range("F63:F67").NumberFormat = $ws.Range("F61").NumberFormat
$ws.Range("G63:G67").NumberFormat = $ws.Range("G61").NumberFormat

# --- Update the sheet view to match where the user ended up scrolled to ---
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("D71").Select() | Out-Null
